$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions refresh)
# Force the Price column cells we touch to Text format so strings like "4.20" / "0.820"
# and multi-dot values like "27.317.05" are preserved exactly (not coerced to Number).
$dPriceCells = @("D2","D3","D5","D6","D10","D11","D12","D13","D14","D15","D16","D17","D19","D21","D22","D24","D25","D27","D29","D30","D32","D33","D35","D38","D39","D41","D42","D45","D46","D47","D49","D50")
foreach ($addr in $dPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.317.05"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "1.659.77"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "219.86"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "0.507"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "20.04"
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.890.00"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "1.651.04"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "4.20"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "67.27"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("D17").Value = "27.295.18"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "223.29"
$ws.Range("E19").Value = "  +5.62%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "4.45"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").Value = "6.73"
$ws.Range("E22").Value = "  +8.72%  "
$ws.Range("E23").Value = "  +4.11%  "
$ws.Range("D24").Value = "9.28"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "147.04"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "7.46"
$ws.Range("E27").Value = "  +4.86%  "
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").Value = "16.08"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "3.42"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").Value = "1.265.49"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("D38").Value = "0.539"
$ws.Range("E38").Value = "  +1.35%  "
$ws.Range("D39").Value = "0.838"
$ws.Range("E39").Value = "  +3.41%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "0.820"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("D42").Value = "5.38"
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("D45").Value = "61.99"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "92.19"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "0.0982"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").Value = "7.68"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  +0.10%  "
